$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrections per row (Diebold-Mariano summary fix)
$ws.Range("B2").Value = "5/10"
$ws.Range("C2").Value = 128

$ws.Range("B3").Value = "5/10"
$ws.Range("C3").Value = 128

$ws.Range("B6").Value = "0/10"
$ws.Range("C6").Value = 0

$ws.Range("B7").Value = "0/10"
$ws.Range("B8").Value = "0/10"
$ws.Range("B9").Value = "0/10"
$ws.Range("B10").Value = "0/10"
